$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 3, pushing the
# previously existing rows 3-5 down to rows 4-6 (their data is unchanged).
$ws.Rows(3).Insert()

$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44883
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112039
$ws.Range("G3").Value = "Ciboulette"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("N3").Value = '$/docena de atados'
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 633
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = "Hortaliza"
